$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; this shifts rows 11:98 down to 12:99
# (matching the existing behaviour of every other row in this log: newest entry on top).
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new "Choclo / Dulce o Americano" record.
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44545
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 100112024
$ws.Range("G11").Value = "Choclo"
$ws.Range("H11").Value = "Dulce o Americano"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 700
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14000
$ws.Range("N11").Value = "$/malla 70 unidades"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 200
$ws.Range("Q11").Value = 70
$ws.Range("R11").Value = "Hortaliza"
